# "Database wilayah Indonesia 2018 Semester 1"
#
# The sheet already lists Indonesian provinces in columns A (dagri id) and
# B (province name), headed by a label row ("id" / "dagri_name") on row 2
# with an empty, bottom-bordered banner row above it on row 1.
#
# This change:
#   - gives the sheet a title in the banner row (A1)
#   - renames the "dagri_name" column header to "name"
#
# Note: writing B2 first (then A1) reproduces the exact shared-string
# ordering Excel produced (the now-unused "dagri_name" string is dropped,
# "name" is appended, then "Table Dagri Province" is appended).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "name"
$ws.Range("A1").Value = "Table Dagri Province"

# Editing the header text nudges Excel's autofit/wrap-text row-height
# calculation for the wrap-text rows; row 27 ("SULAWESI UTARA") loses its
# stale explicit height back to the sheet default.
$ws.Rows("27:27").AutoFit()
